$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue 'D2' '26.843.06'
Set-TextValue 'E2' '  -1.59%  '
Set-TextValue 'D3' '1.811.18'
Set-TextValue 'E3' '  -0.39%  '
Set-TextValue 'E4' '  -0.02%  '
Set-TextValue 'D5' '310.02'
Set-TextValue 'E5' '  -1.05%  '
Set-TextValue 'E6' '  +0.05%  '
Set-TextValue 'D7' '0.4630'
Set-TextValue 'E7' '  -0.28%  '
Set-TextValue 'D8' '0.3694'
Set-TextValue 'E8' '  -1.84%  '
Set-TextValue 'D9' '0.07349'
Set-TextValue 'E9' '  -0.93%  '
Set-TextValue 'D10' '0.8686'
Set-TextValue 'E10' '  -0.05%  '
Set-TextValue 'D11' '20.37'
Set-TextValue 'E11' '  -1.11%  '
Set-TextValue 'D12' '1.842.13'
Set-TextValue 'E12' '  +1.19%  '
Set-TextValue 'D13' '5.338'
Set-TextValue 'E13' '  -1.04%  '
Set-TextValue 'D14' '0.07073'
Set-TextValue 'E14' '  -0.37%  '
Set-TextValue 'E15' '  -2.22%  '
Set-TextValue 'D16' '91.64'
Set-TextValue 'E16' '  -0.41%  '
Set-TextValue 'D17' '1.002'
Set-TextValue 'E17' '  +0.01%  '
Set-TextValue 'D18' '0.000008697'
Set-TextValue 'E18' '  -0.67%  '
Set-TextValue 'E19' '  +0.06%  '
Set-TextValue 'E20' '  -1.85%  '
Set-TextValue 'D21' '26.905.75'
Set-TextValue 'E21' '  -1.47%  '
Set-TextValue 'D22' '5.335'
Set-TextValue 'E22' '  +0.60%  '
Set-TextValue 'D23' '10.54'
Set-TextValue 'E23' '  -3.22%  '
Set-TextValue 'D24' '2.110.75'
Set-TextValue 'E24' '  +2.80%  '
Set-TextValue 'D25' '1.898'
Set-TextValue 'E25' '  -2.28%  '
Set-TextValue 'D26' '151.88'
Set-TextValue 'E27' '  -0.92%  '
Set-TextValue 'D28' '2.112'
Set-TextValue 'E28' '  -6.60%  '
Set-TextValue 'D29' '5.288'
Set-TextValue 'E29' '  +0.13%  '
Set-TextValue 'D30' '115.32'
Set-TextValue 'E30' '  -1.37%  '
Set-TextValue 'D31' '0.08888'
Set-TextValue 'E31' '  -0.21%  '
Set-TextValue 'D32' '0.7555'
Set-TextValue 'E32' '  -2.81%  '
Set-TextValue 'B33' 'ARBITRUM'
Set-TextValue 'C33' 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue 'D33' '1.149'
Set-TextValue 'E33' '  -2.67%  '
Set-TextValue 'B34' 'HuobiToken'
Set-TextValue 'C34' 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue 'D34' '2.921'
Set-TextValue 'E34' '  -0.01%  '
Set-TextValue 'D35' '4.456'
Set-TextValue 'E35' '  -1.48%  '
Set-TextValue 'E36' '  +0.13%  '
Set-TextValue 'D37' '1.096'
Set-TextValue 'E37' '  -0.52%  '
Set-TextValue 'D38' '0.01952'
Set-TextValue 'E38' '  -0.85%  '
Set-TextValue 'D39' '0.05252'
Set-TextValue 'E39' '  -0.05%  '
Set-TextValue 'B40' 'TheSandbox'
Set-TextValue 'C40' 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue 'D40' '0.5338'
Set-TextValue 'E40' '  +0.95%  '
Set-TextValue 'B41' 'MXToken'
Set-TextValue 'C41' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 'D41' '2.912'
Set-TextValue 'E41' '  +0.01%  '
Set-TextValue 'D42' '7.172'
Set-TextValue 'E42' '  -1.48%  '
Set-TextValue 'E43' '  -1.06%  '
Set-TextValue 'E44' '  -1.75%  '
Set-TextValue 'D45' '8.410'
Set-TextValue 'E45' '  -2.04%  '
Set-TextValue 'D46' '0.4930'
Set-TextValue 'E46' '  -2.16%  '
Set-TextValue 'D47' '10.35'
Set-TextValue 'E47' '  -1.56%  '
Set-TextValue 'E48' '  +0.10%  '
Set-TextValue 'D49' '1.670'
Set-TextValue 'E49' '  -0.21%  '
Set-TextValue 'D50' '102.93'
Set-TextValue 'E50' '  -2.40%  '
Set-TextValue 'E51' '  -1.01%  '
